# feat: add 2022-Q1 data
#
# - Renames the existing "总计" sheet to "2022-Q1" and repopulates it with
#   the per-fund holdings table (same shape as the 2020-Q4..2021-Q3 sheets).
# - Adds a brand-new "总计" sheet (placed right after "2022-Q1", i.e. at the
#   end) with the historical summary table plus the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# Reference cells (still on an untouched sheet) that already carry the
# bold/centered/thin-bordered "header" style (s="2") so we can clone the
# exact formatting instead of re-building it from scratch.
$refSheet = $wb.Worksheets.Item("2021-Q3")
$refHeader = $refSheet.Range("B1")
$refA = $refSheet.Range("A2")

function Set-TextValue($range, [string]$text) {
    # Force literal text even for digit-only strings (fund codes, "91.10",
    # "0.1810", ...) so Excel doesn't "helpfully" renormalize them into
    # numbers and eat meaningful leading/trailing zeros.
    $range.Value = "'" + $text
}

function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) "总计" -> "2022-Q1", repopulated with the fund holdings table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q1.Cells.Item(1, $col).Value = $headers[$col - 2]
}
Copy-Format $refHeader $q1.Range("B1:H1")

$funds = @(
    @("005613", "上投摩根富时发达市场REITs指数QDII人民币份额", "4.84", "91.10", "3.74", "0.1810", 4),
    @("005614", "上投摩根富时发达市场REITs指数QDII美钞",       "4.84", "91.10", "3.74", "0.1810", 4),
    @("005615", "上投摩根富时发达市场REITs指数QDII美汇",       "4.84", "91.10", "3.74", "0.1810", 4),
    @("000179", "广发美国房地产指数QDII-人民币",               "2.37", "92.38", "3.03", "0.0718", 6),
    @("000180", "广发美国房地产指数QDII - 美元",                "2.37", "92.38", "3.03", "0.0718", 6),
    @("160140", "南方道琼斯美国精选REIT指数(QDII-LOF)A",        "1.35", "89.10", "3.27", "0.0441", 5),
    @("070031", "嘉实全球房地产(QDII)",                         "0.60", "95.08", "2.92", "0.0175", 7),
    @("320017", "诺安全球收益不动产(QDII)",                     "0.29", "93.32", "5.47", "0.0159", 7),
    @("160141", "南方道琼斯美国精选REIT指数(QDII-LOF)C",        "0.44", "89.10", "3.27", "0.0144", 5)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $rec = $funds[$i]

    $q1.Cells.Item($row, 1).Value = $i
    Set-TextValue $q1.Cells.Item($row, 2) $rec[0]
    $q1.Cells.Item($row, 3).Value = $rec[1]
    Set-TextValue $q1.Cells.Item($row, 4) $rec[2]
    Set-TextValue $q1.Cells.Item($row, 5) $rec[3]
    Set-TextValue $q1.Cells.Item($row, 6) $rec[4]
    Set-TextValue $q1.Cells.Item($row, 7) $rec[5]
    $q1.Cells.Item($row, 8).Value = $rec[6]
}
Copy-Format $refA $q1.Range("A2:A10")

# ---------------------------------------------------------------------
# 2) Brand-new "总计" sheet right after "2022-Q1", with the historical
#    summary table plus the new 2022-Q1 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $total.Cells.Item(1, $col).Value = $totalHeaders[$col - 2]
}
Copy-Format $refHeader $total.Range("B1:D1")

$history = @(
    @("2022-Q1", 9, 0.78),
    @("2021-Q3", 11, 1.07),
    @("2021-Q2", 10, 1.06),
    @("2021-Q1", 12, 0.9),
    @("2020-Q4", 11, 0.75)
)

for ($i = 0; $i -lt $history.Length; $i++) {
    $row = $i + 2
    $rec = $history[$i]

    $total.Cells.Item($row, 1).Value = $i
    $total.Cells.Item($row, 2).Value = $rec[0]
    $total.Cells.Item($row, 3).Value = $rec[1]
    $total.Cells.Item($row, 4).Value = $rec[2]
}
Copy-Format $refA $total.Range("A2:A6")

Write-Host "2022-Q1 + 总计 rebuilt"
